$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 new rows before row 385; existing rows 385-444 shift down to 389-448.
$ws.Rows("385:388").Insert()

# Row 385
$ws.Range("A385").Value = 5
$ws.Range("B385").Value = "Macroferia Regional de Talca"
$ws.Range("C385").Value = "Maule"
$ws.Range("D385").Value = 45077
$ws.Range("E385").Value = 7
$ws.Range("F385").Value = "Fruta"
$ws.Range("G385").Value = 100102
$ws.Range("H385").Value = "Cítricos"
$ws.Range("I385").Value = 100102004
$ws.Range("J385").Value = "Mandarina"
$ws.Range("K385").Value = "Clementina"
$ws.Range("L385").Value = "Primera"
$ws.Range("M385").Value = 350
$ws.Range("N385").Value = 7000
$ws.Range("O385").Value = 8000
$ws.Range("P385").Value = 7286
$ws.Range("Q385").Value = "`$/bandeja 10 kilos"
$ws.Range("R385").Value = "Provincia de Limarí"
$ws.Range("S385").Value = 729
$ws.Range("T385").Value = 10

# Row 386
$ws.Range("A386").Value = 5
$ws.Range("B386").Value = "Macroferia Regional de Talca"
$ws.Range("C386").Value = "Maule"
$ws.Range("D386").Value = 45077
$ws.Range("E386").Value = 7
$ws.Range("F386").Value = "Fruta"
$ws.Range("G386").Value = 100102
$ws.Range("H386").Value = "Cítricos"
$ws.Range("I386").Value = 100102004
$ws.Range("J386").Value = "Mandarina"
$ws.Range("K386").Value = "Clementina"
$ws.Range("L386").Value = "Primera"
$ws.Range("M386").Value = 180
$ws.Range("N386").Value = 10000
$ws.Range("O386").Value = 10000
$ws.Range("P386").Value = 10000
$ws.Range("Q386").Value = "`$/bandeja 18 kilos"
$ws.Range("R386").Value = "Región de O'Higgins"
$ws.Range("S386").Value = 556
$ws.Range("T386").Value = 18

# Row 387
$ws.Range("A387").Value = 5
$ws.Range("B387").Value = "Macroferia Regional de Talca"
$ws.Range("C387").Value = "Maule"
$ws.Range("D387").Value = 45077
$ws.Range("E387").Value = 7
$ws.Range("F387").Value = "Fruta"
$ws.Range("G387").Value = 100102
$ws.Range("H387").Value = "Cítricos"
$ws.Range("I387").Value = 100102004
$ws.Range("J387").Value = "Mandarina"
$ws.Range("K387").Value = "Clementina"
$ws.Range("L387").Value = "Segunda"
$ws.Range("M387").Value = 200
$ws.Range("N387").Value = 8000
$ws.Range("O387").Value = 8000
$ws.Range("P387").Value = 8000
$ws.Range("Q387").Value = "`$/bandeja 18 kilos"
$ws.Range("R387").Value = "Región de O'Higgins"
$ws.Range("S387").Value = 444
$ws.Range("T387").Value = 18

# Row 388
$ws.Range("A388").Value = 5
$ws.Range("B388").Value = "Macroferia Regional de Talca"
$ws.Range("C388").Value = "Maule"
$ws.Range("D388").Value = 45077
$ws.Range("E388").Value = 7
$ws.Range("F388").Value = "Fruta"
$ws.Range("G388").Value = 100102
$ws.Range("H388").Value = "Cítricos"
$ws.Range("I388").Value = 100102004
$ws.Range("J388").Value = "Mandarina"
$ws.Range("K388").Value = "Clementina"
$ws.Range("L388").Value = "Tercera"
$ws.Range("M388").Value = 150
$ws.Range("N388").Value = 6000
$ws.Range("O388").Value = 6000
$ws.Range("P388").Value = 6000
$ws.Range("Q388").Value = "`$/bandeja 18 kilos"
$ws.Range("R388").Value = "Región de O'Higgins"
$ws.Range("S388").Value = 333
$ws.Range("T388").Value = 18
